# Add a new "2022-Q1" sheet (same layout as the existing quarterly sheets)
# right before the "总计" (totals) sheet, and add a matching summary row
# to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4"
#    (i.e. right before "总计"), using "2021-Q4" as a formatting template.
# ---------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$newSheet = $wb.Worksheets.Add($null, $template)
$newSheet.Name = "2022-Q1"

# Copy the template's cell formatting (borders/fonts/alignment for the
# header row and the index column) onto the same range of the new sheet.
# The template sheet only has 4 data rows, so its style is copied in two
# passes to cover all 6 rows needed here without leaving unstyled cells.
$template.Range("A1:H4").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)
$template.Range("A3:H4").Copy()
$newSheet.Range("A5:H6").PasteSpecial(-4122)

# Header row
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

# Keep the fund-code column and the numeric-looking text columns as text
# so values like "002446" keep their leading zero and "12.46" is not
# silently turned into a number.
$newSheet.Range("B2:G6").NumberFormat = "@"

# Row 2
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "002446"
$newSheet.Range("C2").Value = "广发利鑫灵活配置混合A"
$newSheet.Range("D2").Value = "12.46"
$newSheet.Range("E2").Value = "74.35"
$newSheet.Range("F2").Value = "4.25"
$newSheet.Range("G2").Value = "0.5296"
$newSheet.Range("H2").Value = 6

# Row 3
$newSheet.Range("A3").Value = 1
$newSheet.Range("B3").Value = "161605"
$newSheet.Range("C3").Value = "融通蓝筹成长混合"
$newSheet.Range("D3").Value = "4.82"
$newSheet.Range("E3").Value = "71.70"
$newSheet.Range("F3").Value = "4.89"
$newSheet.Range("G3").Value = "0.2357"
$newSheet.Range("H3").Value = 4

# Row 4
$newSheet.Range("A4").Value = 2
$newSheet.Range("B4").Value = "000717"
$newSheet.Range("C4").Value = "融通转型三动力灵活配置混合A"
$newSheet.Range("D4").Value = "3.83"
$newSheet.Range("E4").Value = "94.89"
$newSheet.Range("F4").Value = "6.14"
$newSheet.Range("G4").Value = "0.2352"
$newSheet.Range("H4").Value = 4

# Row 5
$newSheet.Range("A5").Value = 3
$newSheet.Range("B5").Value = "011172"
$newSheet.Range("C5").Value = "广发利鑫灵活配置混合C"
$newSheet.Range("D5").Value = "1.10"
$newSheet.Range("E5").Value = "74.35"
$newSheet.Range("F5").Value = "4.25"
$newSheet.Range("G5").Value = "0.0468"
$newSheet.Range("H5").Value = 6

# Row 6
$newSheet.Range("A6").Value = 4
$newSheet.Range("B6").Value = "009828"
$newSheet.Range("C6").Value = "融通转型三动力灵活配置混合C"
$newSheet.Range("D6").Value = "0.59"
$newSheet.Range("E6").Value = "94.89"
$newSheet.Range("F6").Value = "6.14"
$newSheet.Range("G6").Value = "0.0362"
$newSheet.Range("H6").Value = 4

# ---------------------------------------------------------------------
# 2. Add a "2022-Q1" row at the top of the "总计" (totals) sheet's data,
#    pushing the existing rows down by one.
# ---------------------------------------------------------------------
$totals = $wb.Worksheets.Item("总计")

$totals.Rows.Item(2).Insert()

# Match formatting of the row above/below it (index column style).
$totals.Range("A3").Copy()
$totals.Range("A2").PasteSpecial(-4122)
$totals.Range("B2:D2").ClearFormats()

$totals.Range("A2").Value = 0
$totals.Range("B2").Value = "2022-Q1"
$totals.Range("C2").Value = 5
$totals.Range("D2").Value = 1.08

# Renumber the index column for the rows that got pushed down.
$totals.Range("A3").Value = 1
$totals.Range("A4").Value = 2
$totals.Range("A5").Value = 3
$totals.Range("A6").Value = 4
$totals.Range("A7").Value = 5

$excel.CutCopyMode = $false
